$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.138.57'
$ws.Range('E2').Value = '  +0.42%  '
$ws.Range('D3').Value = '1.678.19'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'214.28"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.68%  '
$ws.Range('D6').Value = "'0.519"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.16%  '
$ws.Range('D8').Value = "'22.71"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.48%  '
$ws.Range('E9').Value = '  +2.10%  '
$ws.Range('E10').Value = '  -0.55%  '
$ws.Range('E11').Value = '  +0.18%  '
$ws.Range('D12').Value = '1.915.14'
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').Value = '1.680.45'
$ws.Range('E13').Value = '  +0.72%  '
$ws.Range('E14').Value = '  +2.25%  '
$ws.Range('D15').Value = "'0.549"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.98%  '
$ws.Range('D16').Value = "'66.50"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('D17').Value = '27.101.85'
$ws.Range('E17').Value = '  +0.31%  '
$ws.Range('D18').Value = "'234.93"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.35%  '
$ws.Range('D19').Value = "'7.87"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.31%  '
$ws.Range('D20').Value = '0.0₃0739'
$ws.Range('E20').Value = '  +0.27%  '
$ws.Range('E22').Value = '  +1.55%  '
$ws.Range('E23').Value = '  +2.61%  '
$ws.Range('E24').Value = '  -1.37%  '
$ws.Range('D25').Value = "'147.78"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.76%  '
$ws.Range('D26').Value = "'7.42"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.36%  '
$ws.Range('D27').Value = "'16.33"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.59%  '
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('E30').Value = '  +0.71%  '
$ws.Range('E31').Value = '  -0.29%  '
$ws.Range('D32').Value = "'3.35"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.15%  '
$ws.Range('D33').Value = '1.542.68'
$ws.Range('E33').Value = '  +0.24%  '
$ws.Range('E34').Value = '  +1.37%  '
$ws.Range('E35').Value = '  -3.18%  '
$ws.Range('D36').Value = "'0.605"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.87%  '
$ws.Range('D37').Value = "'0.939"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.23%  '
$ws.Range('E38').Value = '  -0.31%  '
$ws.Range('E39').Value = '  -1.41%  '
$ws.Range('E40').Value = '  +2.71%  '
$ws.Range('E41').Value = '  +3.85%  '
$ws.Range('D42').Value = "'69.34"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.10%  '
$ws.Range('E43').Value = '  +0.12%  '
$ws.Range('D45').Value = '1.822.44'
$ws.Range('E45').Value = '  +0.16%  '
$ws.Range('D46').Value = "'0.779"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.41%  '
$ws.Range('D47').Value = "'89.94"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.56%  '
$ws.Range('E48').Value = '  +5.95%  '
$ws.Range('E49').Value = '  +2.54%  '
$ws.Range('E50').Value = '  +2.12%  '
$ws.Range('E51').Value = '  -0.29%  '
